# Add new column 'Event' to Card1 by admin
#
# The "Event" / "Correction" / "Serviced by " headers (M1:O1) already exist.
# Making room for data under the new "Event" column means the existing
# M:N data (rows 2-12) needs to slide one column to the right, into N:O,
# leaving column M blank for the new Event values. The trailing space in
# the "Serviced by " header is also tidied up while we're at it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# Tidy the trailing-space typo in the "Serviced by " header.
$ws.Range("O1").Value = "Serviced by"

# Shift the M:N data for every data row one column to the right (M->N,
# N->O), clearing out the new M column so it's ready for Event entries.
$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    $oldM = $ws.Cells.Item($r, 13).Value2
    $oldN = $ws.Cells.Item($r, 14).Value2

    $ws.Cells.Item($r, 15).Value = $oldN
    $ws.Cells.Item($r, 14).Value = $oldM
    $ws.Cells.Item($r, 13).Value = ""
    $ws.Cells.Item($r, 13).Font.Bold = $false
}
